$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ACE_landing_page_data")
$ws2 = $wb.Worksheets.Item("ANSP")

# New data rows for ACE_landing_page_data (year rolled forward, values recalculated from source)
$data = @(
    @(2024, 469.4605934329557, 10467653125.780624, 22297192.293042853, 0.98081881913019253, 152.18400188546104, 314.30043651384142, -0.013990770346900772, 0.040218364367912596, 0.054978323817399977, 0.047861539125722885, 0.038514105249225938, -0.016474730305083396, 100.98603190277626, 102.09317864913423),
    @(2023, 476.12190567234626, 10062938210.229813, 21135213.671842363, 0.93601948588411121, 146.54013952842695, 319.56518678095119, -0.072964061904103206, 0.029726379638513656, 0.11077288088047554, 0.088513539388441087, 0.021614143908957439, -0.078497934256253843, 97.081569949152268, 96.772773756823597),
    @(2022, 513.59595254773615, 9772438979.1416397, 19027484.408053901, 0.85990614908657403, 143.43981081519371, 346.78727119621755, -0.34491825017192079, 0.036839780409229528, 0.58276395378338597, 0.42197440293200161, -0.033503681179227307, -0.35612961410659494, 94.278996701272078, 87.122016951039342),
    @(2021, 784.01810565250082, 9425216088.144846, 12021681.668053685, 0.6047268834892624, 148.41216466318815, 538.59795200090059, -0.25272924176920852, -0.049774400373029382, 0.27159478563925998, 0.25213414892893815, -0.082604403998122189, -0.24587929212811932, 90.929185475562264, 55.044226110144713),
    @(2020, 1049.1754120135931, 9918924613.1075573, 9454019.3179622926, 0.48295694515363163, 161.77553643159669, 714.20655390940976, 1.210617188802134, -0.043077922640073374, -0.56712447446477454, -0.50708837746858282, 0.083696161698424421, 1.2163174140894872, 95.692207735992667, 43.287552553522552),
    @(2019, 474.60746135883858, 10365446516.265041, 21840041.213401809, 0.97980433626891994, 149.2812673416262, 322.2492181711354, -0.0023237134232800827, 0.014159414678369009, 0.016521519378000837, 0.010339197903164399, 0.005633255473421972, -0.0012163576101731355, 100, 100)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws1.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}

# Activate sheet 1 tab and move selection to J20
$ws1.Activate()
$ws1.Range("J20").Select()
